$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 8.487299773058886
$ws.Range("F3").Value = 8.17703893189862
$ws.Range("F4").Value = 7.389685286561035
$ws.Range("F5").Value = 7.23768709675885
$ws.Range("F6").Value = 6.343048961041616
$ws.Range("F7").Value = 6.040248810237871
$ws.Range("F8").Value = 5.132481009527025
$ws.Range("F9").Value = 3.39083832478149
$ws.Range("F10").Value = 3.318523148442117
$ws.Range("F11").Value = 2.176093329063497
$ws.Range("F12").Value = 1.474100674863718
$ws.Range("F13").Value = 0.05652031328557322
